# Auto-generated update script for cryptos.xlsx
# Applies Price (D) and Volume(1h) (E) updates for rows 2-51,
# plus a swap of the WEMIXToken/Decentraland rows (45/46),
# matching a refreshed data pull from coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.878.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.08%  "

$ws.Range("D3").Value = "'1.884.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.17%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").Value = "'323.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "'0.4673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("D8").Value = "'0.3942"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.10%  "

$ws.Range("D9").Value = "'0.07934"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.9837"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.49%  "

$ws.Range("D11").Value = "'22.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").Value = "'1.894.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.06%  "

$ws.Range("D13").Value = "'5.760"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").Value = "'7.025"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").Value = "'0.06993"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").Value = "'88.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.90%  "

$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "'0.00001013"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("D19").Value = "'17.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").Value = "'28.873.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("D22").Value = "'5.363"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("D23").Value = "'11.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").Value = "'2.120"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").Value = "'2.123.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.45%  "

$ws.Range("D26").Value = "'153.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").Value = "'19.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").Value = "'5.792"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("D29").Value = "'2.009"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "

$ws.Range("D30").Value = "'120.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("D31").Value = "'0.09411"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "

$ws.Range("D32").Value = "'0.9422"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("D33").Value = "'5.326"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").Value = "'1.366"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.73%  "

$ws.Range("D35").Value = "'3.344"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'0.05931"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").Value = "'0.02129"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").Value = "'1.161"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").Value = "'7.913"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.00%  "

$ws.Range("D40").Value = "'0.5750"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.19%  "

$ws.Range("D41").Value = "'0.1803"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").Value = "'10.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.11%  "

$ws.Range("D43").Value = "'0.07313"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.38%  "

$ws.Range("E44").Value = "  +2.27%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5364"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.167"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.78%  "

$ws.Range("D47").Value = "'2.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.84%  "

$ws.Range("D48").Value = "'1.854"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").Value = "'114.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "

$ws.Range("D50").Value = "'2.375"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = "'1.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "

Write-Output "Updated cryptos list"
